$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the header row
$ws.Range("B1").Value = "NOMBRES"
$ws.Range("C1").Value = "APELLIDOS"
$ws.Range("D1").Value = "CORREO"
$ws.Range("E1").Value = "CELULAR"
$ws.Range("F1").Value = "EMPRESA"
$ws.Range("G1").Value = "CARGO"

# Update row 2 data
$ws.Range("D2").Value = "gabrielmasutier@gmail.com"
$ws.Range("F2").Value = "La Data Inquisición"
$ws.Range("G2").Value = "Analista"

# Remove row 3 (Liliana Guerrero entry) entirely
$ws.Rows.Item(3).Delete()
